# Auto-generated: reapply the per-row price/volume refresh described
# in the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.083.47'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '2.468.49'
$ws.Range('E3').Value = '  -2.83%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '582.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.61'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.11%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.514'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Value = '2.469.95'
$ws.Range('E10').Value = '  -2.30%  '
$ws.Range('E11').Value = '  -0.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.96'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.12%  '
$ws.Range('E13').Value = '  -3.89%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.56'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.71%  '
$ws.Range('D15').Value = '2.919.58'
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('D16').Value = '66.955.73'
$ws.Range('E16').Value = '  -0.77%  '
$ws.Range('E17').Value = '  -4.71%  '
$ws.Range('D18').Value = '2.447.38'
$ws.Range('E18').Value = '  -3.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.12'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -5.83%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.58'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '353.04'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.66%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.03'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.97%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  -3.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.25'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -7.38%  '
$ws.Range('E26').Value = '  -7.41%  '
$ws.Range('E27').Value = '  -7.39%  '
$ws.Range('E28').Value = '  -59.03%  '
$ws.Range('D29').Value = '2.587.42'
$ws.Range('E29').Value = '  -2.96%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '518.50'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.67%  '
$ws.Range('D31').Value = '0.0₃0905'
$ws.Range('E31').Value = '  -6.64%  '
$ws.Range('E32').Value = '  -8.59%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.77'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.47%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.23'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.68%  '
$ws.Range('E35').Value = '  +0.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.118'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -8.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '158.30'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.67'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.37'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.34%  '
$ws.Range('E40').Value = '  -6.61%  '
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E42').Value = '  -6.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.80'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.326'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.38'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '38.73'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '140.92'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.46%  '
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.45'
$ws.Range('D48').Style = 'Normal'
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.515'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.97%  '
$ws.Range('E50').Value = '  -12.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.59'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.60%  '
